# Bugfix for voice_data message corruption + added decoded message prints
# En/Decryption currently disabled

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: fix the wording + status of the voice_data root-cause item ---
$ws1.Range("B62").Value = "Root cause voice_data message loss and corruption in duplex mode"
$ws1.Range("C62").Value = "Ongoing"

# --- Add a new Sheet2 (placed right after Sheet1) with scratch calculations ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$aValues = @(166, 59, 59, 555, 42, 555, 52, 150)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 1).Value = $aValues[$i]
    $ws2.Cells.Item($r, 2).Value = 0.5
    $ws2.Cells.Item($r, 3).Formula = "=A" + $r + "*B" + $r
}

# Row 9: A9 itself is a formula
$ws2.Range("A9").Formula = "=13500"
$ws2.Range("B9").Value = 0.5
$ws2.Range("C9").Formula = "=A9*B9"

$ws2.Range("A11").Value = 8388

$ws2.Range("B13").Formula = "=7000*B9"

[void]$ws2.Range("C1:C9").Select()

# --- Restore Sheet1 as the active sheet/selection ---
[void]$ws1.Activate()
[void]$ws1.Range("C63").Select()
